# Append 9 new daily COVID summary rows (2022-04-04 .. 2022-04-14, excluding
# the 2022-04-09/04-10 weekend which the source data skips) to the bottom of
# the single data table on Sheet1, rows 589-597.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(
    @("2022-04-04", "overview", "K02000001", "United Kingdom", 21359681, 143382, 210, 165780),
    @("2022-04-05", "overview", "K02000001", "United Kingdom", 21410305, 50202, 368, 166148),
    @("2022-04-06", "overview", "K02000001", "United Kingdom", 21461556, 51253, 233, 169095),
    @("2022-04-07", "overview", "K02000001", "United Kingdom", 21508546, 47126, 317, 169412),
    @("2022-04-08", "overview", "K02000001", "United Kingdom", 21549830, 41384, 347, 169759),
    @("2022-04-11", "overview", "K02000001", "United Kingdom", 21641004, 91304, 348, 170107),
    @("2022-04-12", "overview", "K02000001", "United Kingdom", 21679280, 37819, 288, 170395),
    @("2022-04-13", "overview", "K02000001", "United Kingdom", 21715116, 35926, 651, 171046),
    @("2022-04-14", "overview", "K02000001", "United Kingdom", 21747638, 32608, 350, 171396)
)

$startRow = 589
for ($i = 0; $i -lt $rows.Length; $i++) {
    $r = $startRow + $i
    $data = $rows[$i]

    # Column A holds the date as plain text ("2022-04-04") in every existing
    # row of this sheet, not a real Excel date serial. A leading apostrophe
    # forces Excel to store the literal string instead of auto-converting it
    # to a date value.
    $ws.Cells.Item($r, 1).Value = "'" + $data[0]
    $ws.Cells.Item($r, 2).Value = $data[1]
    $ws.Cells.Item($r, 3).Value = $data[2]
    $ws.Cells.Item($r, 4).Value = $data[3]
    $ws.Cells.Item($r, 5).Value = $data[4]
    $ws.Cells.Item($r, 6).Value = $data[5]
    $ws.Cells.Item($r, 7).Value = $data[6]
    $ws.Cells.Item($r, 8).Value = $data[7]
}
